# Update the Bmp7-Eng LR-pairs sheet with the refreshed TPM-based values.
# Sending clusters now include ECs (in addition to FAPs and MuSCs), so the
# 8-row table (2 senders x 4 targets) grows to 12 rows (3 senders x 4 targets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp7"
$ws.Range("C2").Value = "Eng"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03443933333333333
$ws.Range("H2").Value = 0.103318
$ws.Range("I2").Value = 0.05823261822459219
$ws.Range("J2").Value = 0.0582326182245922
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 218.721583
$ws.Range("N2").Value = 656.164749
$ws.Range("O2").Value = 0.7793342808141792
$ws.Range("P2").Value = 0.7793342808141792
$ws.Range("Q2").Value = 7.532625504131333
$ws.Range("R2").Value = 67.793629537182
$ws.Range("S2").Value = 0.04538267564398922
$ws.Range("T2").Value = 0.04538267564398923

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("C3").Value = "Eng"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03443933333333333
$ws.Range("H3").Value = 0.103318
$ws.Range("I3").Value = 0.05823261822459219
$ws.Range("J3").Value = 0.0582326182245922
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 46.29469433333333
$ws.Range("N3").Value = 138.884083
$ws.Range("O3").Value = 0.164954193449581
$ws.Range("P3").Value = 0.164954193449581
$ws.Range("Q3").Value = 1.594358409710444
$ws.Range("R3").Value = 14.349225687394
$ws.Range("S3").Value = 0.009605714571694977
$ws.Range("T3").Value = 0.009605714571694977

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Eng"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03443933333333333
$ws.Range("H4").Value = 0.103318
$ws.Range("I4").Value = 0.05823261822459219
$ws.Range("J4").Value = 0.0582326182245922
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.515309999999999
$ws.Range("N4").Value = 25.54593
$ws.Range("O4").Value = 0.03034118948727519
$ws.Range("P4").Value = 0.03034118948727519
$ws.Range("Q4").Value = 0.2932615995266666
$ws.Range("R4").Value = 2.63935439574
$ws.Range("S4").Value = 0.001766846903892506
$ws.Range("T4").Value = 0.001766846903892506

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Eng"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03443933333333333
$ws.Range("H5").Value = 0.103318
$ws.Range("I5").Value = 0.05823261822459219
$ws.Range("J5").Value = 0.0582326182245922
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.120231
$ws.Range("N5").Value = 21.360693
$ws.Range("O5").Value = 0.02537033624896462
$ws.Range("P5").Value = 0.02537033624896462
$ws.Range("Q5").Value = 0.2452160088193333
$ws.Range("R5").Value = 2.206944079374
$ws.Range("S5").Value = 0.001477381105015489
$ws.Range("T5").Value = 0.001477381105015489

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("C6").Value = "Eng"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.51625
$ws.Range("H6").Value = 1.54875
$ws.Range("I6").Value = 0.8729143757654733
$ws.Range("J6").Value = 0.8729143757654734
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 218.721583
$ws.Range("N6").Value = 656.164749
$ws.Range("O6").Value = 0.7793342808141792
$ws.Range("P6").Value = 0.7793342808141792
$ws.Range("Q6").Value = 112.91501722375
$ws.Range("R6").Value = 1016.23515501375
$ws.Range("S6").Value = 0.6802920972495433
$ws.Range("T6").Value = 0.6802920972495434

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp7"
$ws.Range("C7").Value = "Eng"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.51625
$ws.Range("H7").Value = 1.54875
$ws.Range("I7").Value = 0.8729143757654733
$ws.Range("J7").Value = 0.8729143757654734
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 46.29469433333333
$ws.Range("N7").Value = 138.884083
$ws.Range("O7").Value = 0.164954193449581
$ws.Range("P7").Value = 0.164954193449581
$ws.Range("Q7").Value = 23.89963594958333
$ws.Range("R7").Value = 215.09672354625
$ws.Range("S7").Value = 0.1439908868049381
$ws.Range("T7").Value = 0.1439908868049382

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bmp7"
$ws.Range("C8").Value = "Eng"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.51625
$ws.Range("H8").Value = 1.54875
$ws.Range("I8").Value = 0.8729143757654733
$ws.Range("J8").Value = 0.8729143757654734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.515309999999999
$ws.Range("N8").Value = 25.54593
$ws.Range("O8").Value = 0.03034118948727519
$ws.Range("P8").Value = 0.03034118948727519
$ws.Range("Q8").Value = 4.3960287875
$ws.Range("R8").Value = 39.5642590875
$ws.Range("S8").Value = 0.02648526048126676
$ws.Range("T8").Value = 0.02648526048126676

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bmp7"
$ws.Range("C9").Value = "Eng"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.51625
$ws.Range("H9").Value = 1.54875
$ws.Range("I9").Value = 0.8729143757654733
$ws.Range("J9").Value = 0.8729143757654734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.120231
$ws.Range("N9").Value = 21.360693
$ws.Range("O9").Value = 0.02537033624896462
$ws.Range("P9").Value = 0.02537033624896462
$ws.Range("Q9").Value = 3.67581925375
$ws.Range("R9").Value = 33.08237328375
$ws.Range("S9").Value = 0.02214613122972511
$ws.Range("T9").Value = 0.02214613122972511

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Bmp7"
$ws.Range("C10").Value = "Eng"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04072033333333334
$ws.Range("H10").Value = 0.122161
$ws.Range("I10").Value = 0.06885300600993445
$ws.Range("J10").Value = 0.06885300600993445
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 218.721583
$ws.Range("N10").Value = 656.164749
$ws.Range("O10").Value = 0.7793342808141792
$ws.Range("P10").Value = 0.7793342808141792
$ws.Range("Q10").Value = 8.906415766954334
$ws.Range("R10").Value = 80.15774190258901
$ws.Range("S10").Value = 0.05365950792064662
$ws.Range("T10").Value = 0.05365950792064662

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Bmp7"
$ws.Range("C11").Value = "Eng"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04072033333333334
$ws.Range("H11").Value = 0.122161
$ws.Range("I11").Value = 0.06885300600993445
$ws.Range("J11").Value = 0.06885300600993445
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 46.29469433333333
$ws.Range("N11").Value = 138.884083
$ws.Range("O11").Value = 0.164954193449581
$ws.Range("P11").Value = 0.164954193449581
$ws.Range("Q11").Value = 1.885135384818111
$ws.Range("R11").Value = 16.966218463363
$ws.Range("S11").Value = 0.01135759207294789
$ws.Range("T11").Value = 0.01135759207294789

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Bmp7"
$ws.Range("C12").Value = "Eng"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04072033333333334
$ws.Range("H12").Value = 0.122161
$ws.Range("I12").Value = 0.06885300600993445
$ws.Range("J12").Value = 0.06885300600993445
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 8.515309999999999
$ws.Range("N12").Value = 25.54593
$ws.Range("O12").Value = 0.03034118948727519
$ws.Range("P12").Value = 0.03034118948727519
$ws.Range("Q12").Value = 0.3467462616366667
$ws.Range("R12").Value = 3.12071635473
$ws.Range("S12").Value = 0.002089082102115919
$ws.Range("T12").Value = 0.002089082102115919

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Bmp7"
$ws.Range("C13").Value = "Eng"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04072033333333334
$ws.Range("H13").Value = 0.122161
$ws.Range("I13").Value = 0.06885300600993445
$ws.Range("J13").Value = 0.06885300600993445
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 7.120231
$ws.Range("N13").Value = 21.360693
$ws.Range("O13").Value = 0.02537033624896462
$ws.Range("P13").Value = 0.02537033624896462
$ws.Range("Q13").Value = 0.2899381797303334
$ws.Range("R13").Value = 2.609443617573
$ws.Range("S13").Value = 0.001746823914224019
$ws.Range("T13").Value = 0.001746823914224019

